$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 61.06878433333333
$ws.Range("H2").Value = 183.206353
$ws.Range("I2").Value = 0.4308066250287063
$ws.Range("J2").Value = 0.4308066250287063
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 9739.471853004883
$ws.Range("R2").Value = 87655.24667704395
$ws.Range("S2").Value = 0.1285218978925219
$ws.Range("T2").Value = 0.1285218978925219
$ws.Range("G3").Value = 61.06878433333333
$ws.Range("H3").Value = 183.206353
$ws.Range("I3").Value = 0.4308066250287063
$ws.Range("J3").Value = 0.4308066250287063
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 10537.93763817715
$ws.Range("R3").Value = 94841.43874359432
$ws.Range("S3").Value = 0.1390584382369474
$ws.Range("T3").Value = 0.1390584382369474
$ws.Range("G4").Value = 61.06878433333333
$ws.Range("H4").Value = 183.206353
$ws.Range("I4").Value = 0.4308066250287063
$ws.Range("J4").Value = 0.4308066250287063
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 4542.766754409145
$ws.Range("R4").Value = 40884.9007896823
$ws.Range("S4").Value = 0.05994626954844412
$ws.Range("T4").Value = 0.05994626954844413
$ws.Range("G5").Value = 61.06878433333333
$ws.Range("H5").Value = 183.206353
$ws.Range("I5").Value = 0.4308066250287063
$ws.Range("J5").Value = 0.4308066250287063
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 3567.309484637176
$ws.Range("R5").Value = 32105.78536173458
$ws.Range("S5").Value = 0.0470741527112799
$ws.Range("T5").Value = 0.04707415271127991
$ws.Range("G6").Value = 61.06878433333333
$ws.Range("H6").Value = 183.206353
$ws.Range("I6").Value = 0.4308066250287063
$ws.Range("J6").Value = 0.4308066250287063
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 4259.316622970083
$ws.Range("R6").Value = 38333.84960673075
$ws.Range("S6").Value = 0.05620586663951293
$ws.Range("T6").Value = 0.05620586663951294
$ws.Range("I7").Value = 0.02554841368886107
$ws.Range("J7").Value = 0.02554841368886107
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 577.5864194173121
$ws.Range("R7").Value = 5198.277774755808
$ws.Range("S7").Value = 0.007621820150089192
$ws.Range("T7").Value = 0.007621820150089191
$ws.Range("I8").Value = 0.02554841368886107
$ws.Range("J8").Value = 0.02554841368886107
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("S8").Value = 0.008246675655853105
$ws.Range("T8").Value = 0.008246675655853103
$ws.Range("I9").Value = 0.02554841368886107
$ws.Range("J9").Value = 0.02554841368886107
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 269.402738005517
$ws.Range("R9").Value = 2424.624642049653
$ws.Range("S9").Value = 0.003555033754240836
$ws.Range("T9").Value = 0.003555033754240836
$ws.Range("I10").Value = 0.02554841368886107
$ws.Range("J10").Value = 0.02554841368886107
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 211.554542513443
$ws.Range("R10").Value = 1903.990882620987
$ws.Range("S10").Value = 0.002791669992169368
$ws.Range("T10").Value = 0.002791669992169368
$ws.Range("I11").Value = 0.02554841368886107
$ws.Range("J11").Value = 0.02554841368886107
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 252.593105104248
$ws.Range("R11").Value = 2273.337945938232
$ws.Range("S11").Value = 0.003333214136508574
$ws.Range("T11").Value = 0.003333214136508574
$ws.Range("G12").Value = 36.843258
$ws.Range("H12").Value = 110.529774
$ws.Range("I12").Value = 0.2599088848306786
$ws.Range("J12").Value = 0.2599088848306786
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 5875.896797050433
$ws.Range("R12").Value = 52883.07117345389
$ws.Range("S12").Value = 0.07753822995489423
$ws.Range("T12").Value = 0.07753822995489423
$ws.Range("G13").Value = 36.843258
$ws.Range("H13").Value = 110.529774
$ws.Range("I13").Value = 0.2599088848306786
$ws.Range("J13").Value = 0.2599088848306786
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 6357.617225063225
$ws.Range("R13").Value = 57218.55502556903
$ws.Range("S13").Value = 0.08389500418210254
$ws.Range("T13").Value = 0.08389500418210254
$ws.Range("G14").Value = 36.843258
$ws.Range("H14").Value = 110.529774
$ws.Range("I14").Value = 0.2599088848306786
$ws.Range("J14").Value = 0.2599088848306786
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 2740.685431905062
$ws.Range("R14").Value = 24666.16888714556
$ws.Range("S14").Value = 0.03616603636737756
$ws.Range("T14").Value = 0.03616603636737757
$ws.Range("G15").Value = 36.843258
$ws.Range("H15").Value = 110.529774
$ws.Range("I15").Value = 0.2599088848306786
$ws.Range("J15").Value = 0.2599088848306786
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 2152.184706853498
$ws.Range("R15").Value = 19369.66236168148
$ws.Range("S15").Value = 0.02840019123364819
$ws.Range("T15").Value = 0.02840019123364819
$ws.Range("G16").Value = 36.843258
$ws.Range("H16").Value = 110.529774
$ws.Range("I16").Value = 0.2599088848306786
$ws.Range("J16").Value = 0.2599088848306786
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 2569.677830611728
$ws.Range("R16").Value = 23127.10047550555
$ws.Range("S16").Value = 0.03390942309265609
$ws.Range("T16").Value = 0.0339094230926561
$ws.Range("G17").Value = 2.119603
$ws.Range("H17").Value = 6.358808999999999
$ws.Range("I17").Value = 0.01495263127961596
$ws.Range("J17").Value = 0.01495263127961596
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 338.0419961426453
$ws.Range("R17").Value = 3042.377965283807
$ws.Range("S17").Value = 0.004460796187652124
$ws.Range("T17").Value = 0.004460796187652124
$ws.Range("G18").Value = 2.119603
$ws.Range("H18").Value = 6.358808999999999
$ws.Range("I18").Value = 0.01495263127961596
$ws.Range("J18").Value = 0.01495263127961596
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 365.7555079166909
$ws.Range("R18").Value = 3291.799571250218
$ws.Range("S18").Value = 0.004826503197665012
$ws.Range("T18").Value = 0.004826503197665012
$ws.Range("G19").Value = 2.119603
$ws.Range("H19").Value = 6.358808999999999
$ws.Range("I19").Value = 0.01495263127961596
$ws.Range("J19").Value = 0.01495263127961596
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 157.6724040941836
$ws.Range("R19").Value = 1419.051636847653
$ws.Range("S19").Value = 0.002080642248913019
$ws.Range("T19").Value = 0.002080642248913019
$ws.Range("G20").Value = 2.119603
$ws.Range("H20").Value = 6.358808999999999
$ws.Range("I20").Value = 0.01495263127961596
$ws.Range("J20").Value = 0.01495263127961596
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 123.8157917847763
$ws.Range("R20").Value = 1114.342126062987
$ws.Range("S20").Value = 0.001633870993151973
$ws.Range("T20").Value = 0.001633870993151974
$ws.Range("G21").Value = 2.119603
$ws.Range("H21").Value = 6.358808999999999
$ws.Range("I21").Value = 0.01495263127961596
$ws.Range("J21").Value = 0.01495263127961596
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 147.834288672248
$ws.Range("R21").Value = 1330.508598050232
$ws.Range("S21").Value = 0.001950818652233826
$ws.Range("T21").Value = 0.001950818652233826
$ws.Range("G22").Value = 38.101267
$ws.Range("H22").Value = 114.303801
$ws.Range("I22").Value = 0.268783445172138
$ws.Range("J22").Value = 0.268783445172138
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 6076.528648168502
$ws.Range("R22").Value = 54688.75783351651
$ws.Range("S22").Value = 0.08018576430506832
$ws.Range("T22").Value = 0.08018576430506831
$ws.Range("G23").Value = 38.101267
$ws.Range("H23").Value = 114.303801
$ws.Range("I23").Value = 0.268783445172138
$ws.Range("J23").Value = 0.268783445172138
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 6574.697367315699
$ws.Range("R23").Value = 59172.27630584128
$ws.Range("S23").Value = 0.08675958989045988
$ws.Range("T23").Value = 0.08675958989045987
$ws.Range("G24").Value = 38.101267
$ws.Range("H24").Value = 114.303801
$ws.Range("I24").Value = 0.268783445172138
$ws.Range("J24").Value = 0.268783445172138
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 2834.265835123079
$ws.Range("R24").Value = 25508.39251610772
$ws.Range("S24").Value = 0.0374009217090726
$ws.Range("T24").Value = 0.0374009217090726
$ws.Range("G25").Value = 38.101267
$ws.Range("H25").Value = 114.303801
$ws.Range("I25").Value = 0.268783445172138
$ws.Range("J25").Value = 0.268783445172138
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 2225.67081741636
$ws.Range("R25").Value = 20031.03735674724
$ws.Range("S25").Value = 0.02936991264573532
$ws.Range("T25").Value = 0.02936991264573532
$ws.Range("G26").Value = 38.101267
$ws.Range("H26").Value = 114.303801
$ws.Range("I26").Value = 0.268783445172138
$ws.Range("J26").Value = 0.268783445172138
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 2657.419198055672
$ws.Range("R26").Value = 23916.77278250105
$ws.Range("S26").Value = 0.03506725662180189
$ws.Range("T26").Value = 0.03506725662180189

Write-Host "Updated 278 cells"